# Changed location of the input CSV for easier access to it.
# Refresh the fight card data: update names/predictions for the first
# 24 fights and remove the trailing 4 rows that are no longer part of
# the card.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Cory", "Sandhagen", 0),
    @("Rob", "Font", 0),
    @("Jessica", "Andrade", 0),
    @("Tatiana", "Suarez", 0),
    @("Dustin", "Jacoby", 0),
    @("Kennedy", "Nzechukwu", 0),
    @("Diego", "Lopes", 1),
    @("Gavin", "Tucker", 0),
    @("Tanner", "Boser", 0),
    @("Aleksa", "Camur", 1),
    @("Ignacio", "Bahamondes", 0),
    @("Ludovit", "Klein", 1),
    @("Kyler", "Phillips", 0),
    @("Raoni", "Barcelos", 0),
    @("Jeremiah", "Wells", 1),
    @("Carlston", "Harris", 0),
    @("Billy", "Quarantillo", 0),
    @("Damon", "Jackson", 0),
    @("Cody", "Durden", 0),
    @("Jake", "Hadley", 0),
    @("Sean", "Woodson", 1),
    @("Dennis", "Buzukja", 0),
    @("Ode", "Osbourne", 0),
    @("Assu", "Almabayev", 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

# Remove the now-unused trailing rows (old rows 26-29), shrinking the
# used range down to A1:D25.
$ws.Range("A26:D29").Delete()
